# Updated cryptos list on Sun Feb  4 09:55:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.940.04"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.301.75"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.58%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.72%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("E9").Value = "  -3.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.93"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.55%  "

# Row 13
$ws.Range("E13").Value = "  +0.68%  "

# Row 14
$ws.Range("E14").Value = "  -2.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.658.07"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.303.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.71%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.877.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.57%  "

# Row 20
$ws.Range("E20").Value = "  -0.64%  "

# Row 21
$ws.Range("E21").Value = "  -2.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
$ws.Range("E24").Value = "  -1.19%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("E26").Value = "  -0.91%  "

# Row 27
$ws.Range("E27").Value = "  -0.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.47"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.70%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.98%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.72%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.91"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.22%  "

# Row 34
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.03"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.08"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0687"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("E39").Value = "  -1.49%  "

# Row 40
$ws.Range("E40").Value = "  -2.73%  "

# Row 41
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.73"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.017.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$ws.Range("E44").Value = "  -2.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.44"
$ws.Range("D47").ClearFormats()

# Row 48
$ws.Range("E48").Value = "  -1.83%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.74%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.54"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.46%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.522.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.69%  "
